$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DropTable")

# Insert a new row before the old row 3 (dropId 5001), for new dropId 1002
$ws.Rows.Item(3).Insert()

# Fill the new row 3 with dropId 1002 data
$ws.Range("A3").Value = 1002
$ws.Range("B3").Value = "Heart"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 2

# Update row 2 (dropId 1001): dropEnum changes from "nExp, nGold, nHeart, nEquip" to "Exp, Gold, Heart, Gacha"
$ws.Range("B2").Value = "Exp, Gold, Heart, Gacha"

# Update row 4 (previously row 3, dropId 5001): dropEnum changes to "Exp, Gold, Heart, Gacha"
$ws.Range("B4").Value = "Exp, Gold, Heart, Gacha"
